$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range('D2')
$cell.NumberFormat = "@"
$cell.Value = '82.690.09'
$cell.Style = "Normal"
$ws.Range('E2').Value = '  +4.06%  '
$cell = $ws.Range('D3')
$cell.NumberFormat = "@"
$cell.Value = '3.187.41'
$cell.Style = "Normal"
$ws.Range('E3').Value = '  -0.14%  '
$ws.Range('E4').Value = '  -0.10%  '
$cell = $ws.Range('D5')
$cell.NumberFormat = "@"
$cell.Value = '219.55'
$cell.Style = "Normal"
$ws.Range('E5').Value = '  +6.45%  '
$cell = $ws.Range('D6')
$cell.NumberFormat = "@"
$cell.Value = '620.33'
$cell.Style = "Normal"
$ws.Range('E6').Value = '  -2.14%  '
$cell = $ws.Range('D7')
$cell.NumberFormat = "@"
$cell.Value = '0.293'
$cell.Style = "Normal"
$ws.Range('E7').Value = '  +21.49%  '
$cell = $ws.Range('D8')
$cell.NumberFormat = "@"
$cell.Value = '0.998'
$cell.Style = "Normal"
$ws.Range('E8').Value = '  -0.10%  '
$cell = $ws.Range('D9')
$cell.NumberFormat = "@"
$cell.Value = '0.582'
$cell.Style = "Normal"
$ws.Range('E9').Value = '  -2.79%  '
$cell = $ws.Range('D10')
$cell.NumberFormat = "@"
$cell.Value = '3.185.84'
$cell.Style = "Normal"
$ws.Range('E10').Value = '  -0.19%  '
$cell = $ws.Range('D11')
$cell.NumberFormat = "@"
$cell.Value = '0.597'
$cell.Style = "Normal"
$ws.Range('E11').Value = '  +0.89%  '
$cell = $ws.Range('D12')
$cell.NumberFormat = "@"
$cell.Value = '0.0000259'
$cell.Style = "Normal"
$ws.Range('E12').Value = '  +0.82%  '
$ws.Range('E13').Value = '  -0.25%  '
$ws.Range('E14').Value = '  -1.52%  '
$cell = $ws.Range('D15')
$cell.NumberFormat = "@"
$cell.Value = '3.763.37'
$cell.Style = "Normal"
$ws.Range('E15').Value = '  -0.45%  '
$cell = $ws.Range('D16')
$cell.NumberFormat = "@"
$cell.Value = '32.36'
$cell.Style = "Normal"
$ws.Range('E16').Value = '  +1.18%  '
$cell = $ws.Range('D17')
$cell.NumberFormat = "@"
$cell.Value = '82.242.61'
$cell.Style = "Normal"
$ws.Range('E17').Value = '  +3.63%  '
$cell = $ws.Range('D18')
$cell.NumberFormat = "@"
$cell.Value = '3.171.36'
$cell.Style = "Normal"
$ws.Range('E18').Value = '  -0.47%  '
$cell = $ws.Range('D19')
$cell.NumberFormat = "@"
$cell.Value = '3.26'
$cell.Style = "Normal"
$ws.Range('E19').Value = '  +10.43%  '
$cell = $ws.Range('D20')
$cell.NumberFormat = "@"
$cell.Value = '14.07'
$cell.Style = "Normal"
$ws.Range('E20').Value = '  -3.10%  '
$cell = $ws.Range('D21')
$cell.NumberFormat = "@"
$cell.Value = '438.79'
$cell.Style = "Normal"
$ws.Range('E21').Value = '  -0.13%  '
$cell = $ws.Range('D22')
$cell.NumberFormat = "@"
$cell.Value = '8.94'
$cell.Style = "Normal"
$ws.Range('E22').Value = '  -4.55%  '
$cell = $ws.Range('D23')
$cell.NumberFormat = "@"
$cell.Value = '5.15'
$cell.Style = "Normal"
$ws.Range('E23').Value = '  -0.64%  '
$cell = $ws.Range('D24')
$cell.NumberFormat = "@"
$cell.Value = '7.31'
$cell.Style = "Normal"
$ws.Range('E24').Value = '  +4.91%  '
$cell = $ws.Range('D25')
$cell.NumberFormat = "@"
$cell.Value = '5.24'
$cell.Style = "Normal"
$ws.Range('E25').Value = '  +8.94%  '
$cell = $ws.Range('D26')
$cell.NumberFormat = "@"
$cell.Value = '11.95'
$cell.Style = "Normal"
$ws.Range('E26').Value = '  +10.64%  '
$ws.Range('B27').Value = 'Litecoin'
$ws.Range('C27').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$cell = $ws.Range('D27')
$cell.NumberFormat = "@"
$cell.Value = '77.94'
$cell.Style = "Normal"
$ws.Range('E27').Value = '  +1.03%  '
$ws.Range('B28').Value = 'WrappedeETH'
$ws.Range('C28').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$cell = $ws.Range('D28')
$cell.NumberFormat = "@"
$cell.Value = '3.337.57'
$cell.Style = "Normal"
$ws.Range('E28').Value = '  -0.48%  '
$cell = $ws.Range('D29')
$cell.NumberFormat = "@"
$cell.Value = '0.999'
$cell.Style = "Normal"
$ws.Range('E29').Value = '  -0.06%  '
$cell = $ws.Range('D30')
$cell.NumberFormat = "@"
$cell.Value = '0.0000121'
$cell.Style = "Normal"
$ws.Range('E30').Value = '  -1.11%  '
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$cell = $ws.Range('D31')
$cell.NumberFormat = "@"
$cell.Value = '9.08'
$cell.Style = "Normal"
$ws.Range('E31').Value = '  -0.09%  '
$ws.Range('B32').Value = 'Binance-PegBSC-USD'
$ws.Range('C32').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$cell = $ws.Range('D32')
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.Style = "Normal"
$ws.Range('E32').Value = '  +0.31%  '
$cell = $ws.Range('D33')
$cell.NumberFormat = "@"
$cell.Value = '569.81'
$cell.Style = "Normal"
$ws.Range('E33').Value = '  +5.36%  '
$ws.Range('E34').Value = '  -1.87%  '
$ws.Range('E35').Value = '  +23.00%  '
$cell = $ws.Range('D36')
$cell.NumberFormat = "@"
$cell.Value = '0.153'
$cell.Style = "Normal"
$ws.Range('E36').Value = '  -2.28%  '
$cell = $ws.Range('D37')
$cell.NumberFormat = "@"
$cell.Value = '1.99'
$cell.Style = "Normal"
$ws.Range('E37').Value = '  -2.09%  '
$cell = $ws.Range('D38')
$cell.NumberFormat = "@"
$cell.Value = '22.69'
$cell.Style = "Normal"
$ws.Range('E38').Value = '  -1.35%  '
$cell = $ws.Range('D39')
$cell.NumberFormat = "@"
$cell.Value = '6.25'
$cell.Style = "Normal"
$ws.Range('E39').Value = '  +13.00%  '
$cell = $ws.Range('D40')
$cell.NumberFormat = "@"
$cell.Value = '0.998'
$cell.Style = "Normal"
$ws.Range('E40').Value = '  -0.10%  '
$cell = $ws.Range('D41')
$cell.NumberFormat = "@"
$cell.Value = '0.405'
$cell.Style = "Normal"
$ws.Range('E41').Value = '  -0.81%  '
$cell = $ws.Range('D42')
$cell.NumberFormat = "@"
$cell.Value = '20.86'
$cell.Style = "Normal"
$ws.Range('E42').Value = '  +4.25%  '
$ws.Range('E43').Value = '  +11.78%  '
$cell = $ws.Range('D44')
$cell.NumberFormat = "@"
$cell.Value = '2.99'
$cell.Style = "Normal"
$ws.Range('E44').Value = '  +13.17%  '
$cell = $ws.Range('D45')
$cell.NumberFormat = "@"
$cell.Value = '160.40'
$cell.Style = "Normal"
$ws.Range('E45').Value = '  -2.05%  '
$ws.Range('E46').Value = '  +0.03%  '
$cell = $ws.Range('D47')
$cell.NumberFormat = "@"
$cell.Value = '186.37'
$cell.Style = "Normal"
$ws.Range('E47').Value = '  -2.85%  '
$cell = $ws.Range('D48')
$cell.NumberFormat = "@"
$cell.Value = '44.70'
$cell.Style = "Normal"
$ws.Range('E48').Value = '  +3.59%  '
$ws.Range('E49').Value = '  -0.30%  '
$cell = $ws.Range('D50')
$cell.NumberFormat = "@"
$cell.Value = '0.767'
$cell.Style = "Normal"
$ws.Range('E50').Value = '  -3.93%  '
$cell = $ws.Range('D51')
$cell.NumberFormat = "@"
$cell.Value = '25.93'
$cell.Style = "Normal"
$ws.Range('E51').Value = '  +1.16%  '
